$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 249.697716503316

$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 1766.335244827366
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 2016.558538421709
